# showPorPeriodoActividad demo-data refresh for periodoController
# Replaces the old "prueba" test data with a fresh sample period/activity
# and a fresh integrante record on the second sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1. Actividades")
$ws2 = $wb.Worksheets.Item("Hoja2. Integrantes")

# --- Hoja1. Actividades : header block -------------------------------
$ws1.Range("E4").Value = "qqqqqqqqqq"
$ws1.Range("E5").Value = "Prueba Coordinador"
$ws1.Range("E6").Value = "2020-1"

# --- Hoja1. Actividades : month-letter header row (row 9) ------------
$ws1.Range("D9").Value = "F"
$ws1.Range("E9").Value = "M"
$ws1.Range("F9").Value = "A"
$ws1.Range("G9").Value = "M"
$ws1.Range("H9").Value = "J"

# --- Hoja1. Actividades : first activity row (row 10) -----------------
$ws1.Range("A10").Value = "Actividad 1"
$ws1.Range("B10").Value = "Brayan Legarda"
$ws1.Range("C10").Value = "Aula"
$ws1.Range("E10").Value = "X"
$ws1.Range("F10").Value = "X"
$ws1.Range("G10").Value = "X"
$ws1.Range("I10").Value = "FH05"
$ws1.Range("J10").Value = "No se realizó"

# --- Hoja1. Actividades : second activity row (row 11) is now blank --
$ws1.Range("A11:J11").ClearContents()

# --- Hoja2. Integrantes : replace the sample "estudiante" record -----
$ws2.Range("A6").Value = "David Vergara"
$ws2.Range("C6").Value = 147
$ws2.Range("D6").Value = 147
$ws2.Range("E6").Value = "david_vergara82141@elpoli.edu.co"

# --- Hoja2. Integrantes : second record row is now blank --------------
$ws2.Range("A7:E7").ClearContents()
